$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"5.210326836045978"
$ws.Range("B3").Value = [double]"4.106868953675822"
$ws.Range("B4").Value = [double]"63.8819981799782"
$ws.Range("B5").Value = [double]"1.35677920674344"
$ws.Range("B6").Value = [double]"2.271584767054895"
$ws.Range("B7").Value = [double]"0.2349959708703826"
$ws.Range("B8").Value = [double]"13.95314124366148"
$ws.Range("B9").Value = [double]"0.9002255056896056"
$ws.Range("B10").Value = [double]"0.9068768538645399"
$ws.Range("B11").Value = [double]"0.006651348174934291"
$ws.Range("B12").Value = [double]"0.09106548413833929"
$ws.Range("B13").Value = [double]"0.4157981423550783"
$ws.Range("B14").Value = [double]"0.02925434193076306"
$ws.Range("B15").Value = [double]"45768.00322883359"
$ws.Range("B16").Value = [double]"-0.7962081275249386"
$ws.Range("B17").Value = [double]"1.272774854829956"
$ws.Range("B18").Value = [double]"4.28599968156434"
$ws.Range("B19").Value = [double]"3.171861973472509"
$ws.Range("B20").Value = [double]"0.1283109899121771"
$ws.Range("B22").Value = [double]"0.5499408619047962"
$ws.Range("B23").Value = [double]"1.272774854829956"
$ws.Range("B24").Value = [double]"0.3036584685651698"
$ws.Range("B25").Value = [double]"0.07591461714129245"
$ws.Range("B27").Value = [double]"0.02305211638283606"
$ws.Range("B28").Value = [double]"1.272774854829956"
$ws.Range("B29").Value = [double]"0.3291686249539809"
$ws.Range("B30").Value = [double]"0.1645843124769904"
$ws.Range("B32").Value = [double]"0.05417599182704726"
$ws.Range("B36").Value = [double]"-4.045997172966031e-16"
$ws.Range("B37").Value = [double]"-3.306163351522417e-14"
$ws.Range("B38").Value = [double]"0.001229176673333189"
$ws.Range("B39").Value = [double]"1.056428220328309"
$ws.Range("B40").Value = [double]"60.52887838332872"
$ws.Range("B41").Value = [double]"1.055397025125994"
$ws.Range("B42").Value = [double]"-32.28768987698245"
$ws.Range("B43").Value = [double]"1.510875294266446e-11"
$ws.Range("B45").Value = [double]"0"
$ws.Range("B47").Value = [double]"0"
$ws.Range("B48").Value = [double]"1.510875294266446e-11"
$ws.Range("B49").Value = [double]"6.366404315775473e-11"
$ws.Range("B50").Value = [double]"0.03132118918625627"
$ws.Range("B51").Value = [double]"-1.594607311557835e-17"
$ws.Range("B52").Value = [double]"0.05166237627385222"
$ws.Range("B53").Value = [double]"0.238880615335939"
$ws.Range("B54").Value = [double]"0.05375923065809298"
$ws.Range("B55").Value = [double]"0.2911879055335431"
$ws.Range("B56").Value = [double]"-3.766543489042926e-21"
$ws.Range("B57").Value = [double]"-0.003329638992196813"
$ws.Range("B58").Value = [double]"1.461716366217535e-19"
$ws.Range("B59").Value = [double]"0.1448383823332065"
$ws.Range("B60").Value = [double]"0.003164942144377774"
$ws.Range("B61").Value = [double]"0.01612593386515407"
$ws.Range("B66").Value = [double]"0.01862084684119406"
$ws.Range("B68").Value = [double]"3.875673308855906e-05"
$ws.Range("B69").Value = [double]"0.01740710895336168"
$ws.Range("B70").Value = [double]"-0.03207774747804426"
$ws.Range("B73").Value = [double]"-0.03207774747804426"
$ws.Range("B74").Value = [double]"2.139716888864142"
$ws.Range("B75").Value = [double]"0.1185746731050294"
$ws.Range("B76").Value = [double]"1.196860237688664"
$ws.Range("B79").Value = [double]"1.196860237688664"
$ws.Range("B80").Value = [double]"0.1518292342825849"
$ws.Range("B82").Value = [double]"1.108190542352966"
$ws.Range("B85").Value = [double]"1.108190542352966"
$ws.Range("B87").Value = [double]"0.3291686249539809"
$ws.Range("B88").Value = [double]"-0.0007369563660958242"
$ws.Range("B89").Value = [double]"6.05256449787283e-16"
$ws.Range("B90").Value = [double]"5.571860697434591"
$ws.Range("B91").Value = [double]"-6.223703739549037e-17"
$ws.Range("B92").Value = [double]"7.176195884273442e-10"
$ws.Range("B93").Value = [double]"6.046997146245738e-16"
$ws.Range("B94").Value = [double]"0.0007369563660958242"
$ws.Range("B95").Value = [double]"-0.3983050109420319"
$ws.Range("B96").Value = [double]"6.223703739549037e-17"
$ws.Range("B97").Value = [double]"1.877043210696723e-17"
$ws.Range("B98").Value = [double]"5.557606075728959"
$ws.Range("B99").Value = [double]"6.05256449787283e-16"
$ws.Range("B100").Value = [double]"0.3983050109420319"
$ws.Range("B101").Value = [double]"6.223703739549037e-17"
$ws.Range("B102").Value = [double]"7.176195884273442e-10"
$ws.Range("B103").Value = [double]"-6.046997146245738e-16"
$ws.Range("B104").Value = [double]"0.6077654995871973"
$ws.Range("B105").Value = [double]"6.618928790038293e-17"
$ws.Range("B106").Value = [double]"0.04355761107654907"
$ws.Range("B107").Value = [double]"1.587980113542457e-18"
$ws.Range("B108").Value = [double]"6.116158223273728e-10"
$ws.Range("B109").Value = [double]"-1.542893366512013e-17"
$ws.Range("B110").Value = [double]"5.479840509429206"
$ws.Range("B111").Value = [double]"0.323961536191663"
$ws.Range("B112").Value = [double]"0"
$ws.Range("B113").Value = [double]"-0"
$ws.Range("B114").Value = [double]"0.09126579909611049"
$ws.Range("B115").Value = [double]"-0"
$ws.Range("B116").Value = [double]"4.27943880702091"
$ws.Range("B117").Value = [double]"0.6621475790241521"
$ws.Range("B118").Value = [double]"0.07684489244127322"
$ws.Range("B119").Value = [double]"0.009779472354313834"
$ws.Range("B121").Value = [double]"-1.382357769919018e-18"
$ws.Range("B122").Value = [double]"-0.09140827753479036"
$ws.Range("B123").Value = [double]"7.115076756936123e-20"
$ws.Range("B124").Value = [double]"0.3036584685651698"
$ws.Range("B126").Value = [double]"-0"
$ws.Range("B127").Value = [double]"0.0135150811132985"
$ws.Range("B128").Value = [double]"6.05256449787283e-16"
$ws.Range("B129").Value = [double]"6.363665285524811e-17"
$ws.Range("B130").Value = [double]"0.0004501367408964241"
$ws.Range("B131").Value = [double]"-6.048062235074575e-16"
$ws.Range("B132").Value = [double]"0.3291686249539809"
$ws.Range("B134").Value = [double]"0.0009206738584802246"
$ws.Range("B135").Value = [double]"0.002289169690781615"
$ws.Range("B137").Value = [double]"-1.725768983872683e-20"
$ws.Range("B138").Value = [double]"-0.0003076575845969641"
$ws.Range("B139").Value = [double]"3.535811531435405e-20"
$ws.Range("B142").Value = [double]"0.3495452593500569"
$ws.Range("B143").Value = [double]"0.04875975159197503"
$ws.Range("B144").Value = [double]"6.879150848173862"
$ws.Range("B145").Value = [double]"0.311580471133689"
$ws.Range("B146").Value = [double]"-1.059274657506354e-21"
$ws.Range("B147").Value = [double]"1.568159907414711e-16"
$ws.Range("B148").Value = [double]"-0.6264532024496201"
$ws.Range("B149").Value = [double]"2.356526807207588e-18"
$ws.Range("B150").Value = [double]"0.04300589160287824"
$ws.Range("B151").Value = [double]"-4.01549481311323e-06"
$ws.Range("B152").Value = [double]"5.331388882029812e-07"
$ws.Range("B153").Value = [double]"-0.1247735064688257"
$ws.Range("B154").Value = [double]"-0.09436430209381001"
$ws.Range("B155").Value = [double]"4.717204621872959e-07"
$ws.Range("B156").Value = [double]"0.02515180150071441"
$ws.Range("B157").Value = [double]"0.8952915645781278"
$ws.Range("B158").Value = [double]"-0.002097453548732808"
$ws.Range("B159").Value = [double]"-0.002531269012819026"
$ws.Range("B160").Value = [double]"-0.1958710976282254"
$ws.Range("B161").Value = [double]"-1.117638415059647"
$ws.Range("B162").Value = [double]"0.0003366171735603111"
$ws.Range("B163").Value = [double]"-0.1472961503653625"
$ws.Range("B164").Value = [double]"1.207210182020857"
$ws.Range("B165").Value = [double]"0.08159611729089548"
$ws.Range("B166").Value = [double]"-4.867160655724147e-17"
$ws.Range("B167").Value = [double]"-2.881591893142269e-15"
$ws.Range("B168").Value = [double]"-41.78336518027699"
$ws.Range("B169").Value = [double]"7.299010529158055e-18"
$ws.Range("B170").Value = [double]"2.461263404196501e-05"
$ws.Range("B171").Value = [double]"-1.836153689510622e-05"
$ws.Range("B172").Value = [double]"0.09627448682060243"
$ws.Range("B173").Value = [double]"0.1202862992846506"
$ws.Range("B174").Value = [double]"-9.735078960994101e-06"
$ws.Range("B175").Value = [double]"-0.02765187892277923"

Write-Host "Updated 147 cells"